$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 1.14
$ws.Range("K2").Value = 5.5
$ws.Range("AB2").Value = 34
$ws.Range("G3").Value = 1.83
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 5.25
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 1.62
$ws.Range("M3").Value = 2.2
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.4
$ws.Range("J4").Value = 1.17
$ws.Range("L4").Value = 1.67
$ws.Range("H6").Value = 4.35
$ws.Range("I6").Value = 3.9
$ws.Range("N6").Value = 1.25
$ws.Range("O6").Value = 3.2
$ws.Range("S6").Value = 2.9
$ws.Range("T6").Value = 17
$ws.Range("V6").Value = 9.5
$ws.Range("Z6").Value = 30
$ws.Range("AA6").Value = 11.25
$ws.Range("AB6").Value = 11.75
$ws.Range("AD6").Value = 28
$ws.Range("AE6").Value = 37
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 22
$ws.Range("G14").Value = 2.87
$ws.Range("H14").Value = 3.35
$ws.Range("I14").Value = 2.2
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 7.4
$ws.Range("L14").Value = 1.29
$ws.Range("M14").Value = 3.25
$ws.Range("N14").Value = 1.87
$ws.Range("O14").Value = 1.83
$ws.Range("P14").Value = 1.4
$ws.Range("Q14").Value = 2.72
$ws.Range("R14").Value = 1.72
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 9.5
$ws.Range("U14").Value = 15
$ws.Range("X14").Value = 25
$ws.Range("Y14").Value = 32
$ws.Range("Z14").Value = 7.4
$ws.Range("AA14").Value = 6.6
$ws.Range("AB14").Value = 14
$ws.Range("AC14").Value = 60
$ws.Range("AD14").Value = 8
$ws.Range("AH14").Value = 17.5
$ws.Range("AI14").Value = 27
$ws.Range("AJ14").Value = 450
